$d = $word.ActiveDocument

function Replace-AllText($old, $new) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title / H1 heading and the bolded "Play ..." line near the end (both occurrences)
Replace-AllText "Play Legend of Cleopatra Megaways for Free | Review" "Play Legend of Cleopatra Megaways for Free"

# "What we like" bullet list
Replace-AllText "Megaways gameplay system with a high number of ways to win" "Megaways gameplay system offers unlimited ways to win"
Replace-AllText "Tumbling Reels feature for potential multiple consecutive wins" "Tumbling Reels feature can lead to multiple consecutive wins"
Replace-AllText "Free spins feature with applied multiplier" "Free spins feature with increasing multipliers"
Replace-AllText "Immersion in the ancient Egyptian theme" "Exciting and potentially lucrative gameplay experience"

# "What we don't like" bullet list
Replace-AllText "Graphics not particularly groundbreaking" "Graphics are not groundbreaking"

# Meta description (italic) line near the end
Replace-AllText "Read our review of Legend of Cleopatra Megaways and play for free. Features include Megaways gameplay, Tumbling Reels, and free spins with multiplier." "Read our review of Legend of Cleopatra Megaways and play this game for free. Experience unlimited ways to win and exciting gameplay."
